# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.930.24"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "2.906.86"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "370.08"
$ws.Range("E5").Value = "  +4.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.37"
$ws.Range("E6").Value = "  -5.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.540"
$ws.Range("E7").Value = "  -3.21%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.582"
$ws.Range("E9").Value = "  -4.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.72"
$ws.Range("E10").Value = "  -4.05%  "
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0833"
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.25"
$ws.Range("D14").Value = "3.355.73"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("E15").Value = "  -3.24%  "
$ws.Range("D16").Value = "2.893.21"
$ws.Range("E16").Value = "  -2.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.920"
$ws.Range("E17").Value = "  -7.40%  "
$ws.Range("D18").Value = "50.855.33"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.20"
$ws.Range("E19").Value = "  -6.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.16"
$ws.Range("E20").Value = "  -3.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.89"
$ws.Range("E21").Value = "  -4.73%  "
$ws.Range("D22").Value = "0.0₃0940"
$ws.Range("E22").Value = "  -3.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.83"
$ws.Range("E23").Value = "  -2.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "258.09"
$ws.Range("E24").Value = "  -2.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.66"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.166"
$ws.Range("E28").Value = "  -5.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "25.54"
$ws.Range("E29").Value = "  -4.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.09"
$ws.Range("E30").Value = "  -4.30%  "
$ws.Range("E31").Value = "  -4.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.26"
$ws.Range("E32").Value = "  +3.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.85"
$ws.Range("E33").Value = "  -4.31%  "
$ws.Range("E34").Value = "  -2.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.32"
$ws.Range("E35").Value = "  +1.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "34.09"
$ws.Range("E36").Value = "  -5.58%  "
$ws.Range("E37").Value = "  +0.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0418"
$ws.Range("E39").Value = "  -6.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.98"
$ws.Range("E40").Value = "  -4.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.59"
$ws.Range("E41").Value = "  -4.55%  "
$ws.Range("E42").Value = "  -6.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.112"
$ws.Range("E43").Value = "  -3.75%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.90"
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "118.97"
$ws.Range("E45").Value = "  -3.38%  "
$ws.Range("E46").Value = "  -1.17%  "
$ws.Range("D47").Value = "2.013.17"
$ws.Range("E47").Value = "  -4.73%  "
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.12"
$ws.Range("E49").Value = "  -6.48%  "
$ws.Range("D50").Value = "3.191.12"
$ws.Range("E50").Value = "  -1.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.236"
$ws.Range("E51").Value = "  -0.64%  "

Write-Host "Applied 90 cell updates"
